$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.503243207931519
$ws.Range("B1").Value = 2.62055516242981
$ws.Range("C1").Value = 2.205363750457764
$ws.Range("D1").Value = 2.292547225952148
$ws.Range("E1").Value = 2.619365215301514
